$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-02-14 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-02-15 Saturday", 2) | Out-Null
$d.Content.Find.Execute("79÷2=39, 1", $true, $false, $false, $false, $false, $true, 1, $false, "55÷4=13, 3", 2) | Out-Null
$d.Content.Find.Execute("51÷2=25, 1", $true, $false, $false, $false, $false, $true, 1, $false, "77÷8=9, 5", 2) | Out-Null
$d.Content.Find.Execute("51÷6=8, 3", $true, $false, $false, $false, $false, $true, 1, $false, "99÷6=16, 3", 2) | Out-Null
$d.Content.Find.Execute("85÷5=17, 0", $true, $false, $false, $false, $false, $true, 1, $false, "30÷2=15, 0", 2) | Out-Null
$d.Content.Find.Execute("28÷5=5, 3", $true, $false, $false, $false, $false, $true, 1, $false, "53÷3=17, 2", 2) | Out-Null
$d.Content.Find.Execute("46÷5=9, 1", $true, $false, $false, $false, $false, $true, 1, $false, "58÷7=8, 2", 2) | Out-Null
$d.Content.Find.Execute("86÷7=12, 2", $true, $false, $false, $false, $false, $true, 1, $false, "24÷9=2, 6", 2) | Out-Null
$d.Content.Find.Execute("28÷3=9, 1", $true, $false, $false, $false, $false, $true, 1, $false, "78÷3=26, 0", 2) | Out-Null
$d.Content.Find.Execute("83÷3=27, 2", $true, $false, $false, $false, $false, $true, 1, $false, "63÷9=7, 0", 2) | Out-Null
$d.Content.Find.Execute("39÷5=7, 4", $true, $false, $false, $false, $false, $true, 1, $false, "12÷3=4, 0", 2) | Out-Null
$d.Content.Find.Execute("60÷5=12, 0", $true, $false, $false, $false, $false, $true, 1, $false, "67÷2=33, 1", 2) | Out-Null
$d.Content.Find.Execute("10÷9=1, 1", $true, $false, $false, $false, $false, $true, 1, $false, "64÷3=21, 1", 2) | Out-Null
$d.Content.Find.Execute("82÷9=9, 1", $true, $false, $false, $false, $false, $true, 1, $false, "59÷6=9, 5", 2) | Out-Null
$d.Content.Find.Execute("14÷2=7, 0", $true, $false, $false, $false, $false, $true, 1, $false, "34÷8=4, 2", 2) | Out-Null
$d.Content.Find.Execute("39÷8=4, 7", $true, $false, $false, $false, $false, $true, 1, $false, "66÷5=13, 1", 2) | Out-Null
$d.Content.Find.Execute("49÷8=6, 1", $true, $false, $false, $false, $false, $true, 1, $false, "96÷9=10, 6", 2) | Out-Null
$d.Content.Find.Execute("32÷3=10, 2", $true, $false, $false, $false, $false, $true, 1, $false, "44÷4=11, 0", 2) | Out-Null
$d.Content.Find.Execute("20÷4=5, 0", $true, $false, $false, $false, $false, $true, 1, $false, "59÷2=29, 1", 2) | Out-Null
$d.Content.Find.Execute("73÷9=8, 1", $true, $false, $false, $false, $false, $true, 1, $false, "10÷8=1, 2", 2) | Out-Null
$d.Content.Find.Execute("77÷9=8, 5", $true, $false, $false, $false, $false, $true, 1, $false, "12÷6=2, 0", 2) | Out-Null
$d.Content.Find.Execute("17÷3=5, 2", $true, $false, $false, $false, $false, $true, 1, $false, "32÷8=4, 0", 2) | Out-Null
$d.Content.Find.Execute("94÷2=47, 0", $true, $false, $false, $false, $false, $true, 1, $false, "22÷9=2, 4", 2) | Out-Null
$d.Content.Find.Execute("51÷8=6, 3", $true, $false, $false, $false, $false, $true, 1, $false, "39÷8=4, 7", 2) | Out-Null
$d.Content.Find.Execute("64÷2=32, 0", $true, $false, $false, $false, $false, $true, 1, $false, "55÷4=13, 3", 2) | Out-Null
$d.Content.Find.Execute("59÷5=11, 4", $true, $false, $false, $false, $false, $true, 1, $false, "27÷5=5, 2", 2) | Out-Null
